$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "minicrewmate"
$ws.Range("M18").Value = "ミニクルーメイト"
$ws.Range("A19").Value = "ligftdown"
$ws.Range("M19").Value = "ライトダウン"
$ws.Range("A21").Value = "CAT"
$ws.Range("M21").Value = "猫(青)"
$ws.Range("A22").Value = "Gentle"
$ws.Range("M22").Value = "紳士"

$ws.Range("L19").Select()
